$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginDetails")

$ws.Range("C1").Value = "googleEmail"
$ws.Range("C2").Value = "ritu@csn.com"

$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:ritu@csn.com")

$ws.Activate()
$ws.Range("C2").Select()
